$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 779-781, pushing the existing rows 779-836 down to 782-839
$ws.Range("A779:R781").EntireRow.Insert()

# --- Row 779 ---
$ws.Cells.Item(779,1).Value  = 9
$ws.Cells.Item(779,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(779,3).Value  = "Metropolitana"
$ws.Cells.Item(779,4).Value  = 45265
$ws.Cells.Item(779,5).Value  = 13
$ws.Cells.Item(779,6).Value  = 100112028
$ws.Cells.Item(779,7).Value  = "Sandia"
$ws.Cells.Item(779,8).Value  = "Sin especificar"
$ws.Cells.Item(779,9).Value  = "Primera"
$ws.Cells.Item(779,10).Value = 520
$ws.Cells.Item(779,11).Value = 600
$ws.Cells.Item(779,12).Value = 700
$ws.Cells.Item(779,13).Value = 650
$ws.Cells.Item(779,14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(779,15).Value = "Perú"
$ws.Cells.Item(779,16).Value = 650
$ws.Cells.Item(779,17).Value = 1
$ws.Cells.Item(779,18).Value = "Hortaliza"

# --- Row 780 ---
$ws.Cells.Item(780,1).Value  = 9
$ws.Cells.Item(780,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(780,3).Value  = "Metropolitana"
$ws.Cells.Item(780,4).Value  = 45265
$ws.Cells.Item(780,5).Value  = 13
$ws.Cells.Item(780,6).Value  = 100112028
$ws.Cells.Item(780,7).Value  = "Sandia"
$ws.Cells.Item(780,8).Value  = "Sin especificar"
$ws.Cells.Item(780,9).Value  = "Primera"
$ws.Cells.Item(780,10).Value = 340
$ws.Cells.Item(780,11).Value = 900
$ws.Cells.Item(780,12).Value = 900
$ws.Cells.Item(780,13).Value = 900
$ws.Cells.Item(780,14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(780,15).Value = "Región de O'Higgins"
$ws.Cells.Item(780,16).Value = 900
$ws.Cells.Item(780,17).Value = 1
$ws.Cells.Item(780,18).Value = "Hortaliza"

# --- Row 781 ---
$ws.Cells.Item(781,1).Value  = 9
$ws.Cells.Item(781,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(781,3).Value  = "Metropolitana"
$ws.Cells.Item(781,4).Value  = 45265
$ws.Cells.Item(781,5).Value  = 13
$ws.Cells.Item(781,6).Value  = 100112028
$ws.Cells.Item(781,7).Value  = "Sandia"
$ws.Cells.Item(781,8).Value  = "Sin especificar"
$ws.Cells.Item(781,9).Value  = "Segunda"
$ws.Cells.Item(781,10).Value = 430
$ws.Cells.Item(781,11).Value = 500
$ws.Cells.Item(781,12).Value = 500
$ws.Cells.Item(781,13).Value = 500
$ws.Cells.Item(781,14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(781,15).Value = "Perú"
$ws.Cells.Item(781,16).Value = 500
$ws.Cells.Item(781,17).Value = 1
$ws.Cells.Item(781,18).Value = "Hortaliza"
